$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.356.07"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "3.079.19"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.06"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.25"
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.25"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").Value = "3.610.84"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.54"
$ws.Range("E14").Value = "  -5.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000163"
$ws.Range("E15").Value = "  -3.26%  "
$ws.Range("D16").Value = "57.473.52"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").Value = "3.079.97"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.10"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.10"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.13"
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "333.82"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.501"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("E24").Value = "  +0.95%  "
$ws.Range("E25").Value = "  +3.32%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("D27").Value = "0.0₃0906"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.40"
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.09"
$ws.Range("E29").Value = "  -2.43%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -0.50%  "
$ws.Range("E32").Value = "  -4.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "154.63"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.87"
$ws.Range("E34").Value = "  +10.96%  "
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0674"
$ws.Range("E38").Value = "  -1.19%  "
$ws.Range("D39").Value = "3.121.86"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.75"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.670"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "2.294.78"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0254"
$ws.Range("E45").Value = "  +4.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.38"
$ws.Range("E46").Value = "  -1.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.940"
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.05"
$ws.Range("E48").Value = "  -1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.88"
$ws.Range("E49").Value = "  -3.83%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "250.90"
$ws.Range("E51").Value = "  +8.10%  "
